# Balance the job the monster attr
# Update the "People" sheet: column F (RightMon) rows 4-36 change monster id
# from 11001001 to 11001003, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("People")
$ws.Activate()

$ws.Range("F4:F36").Value = 11001003

$ws.Range("F8").Select()
